# Adds two new columns, I ("I0") and J ("IF"), to Sheet1, mirroring the
# diff that introduced this data alongside the existing columns A-H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row (bold / bordered / centered, matching the rest of row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Reuse the exact header formatting already used by H1 (and the rest of
# row 1) for the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-19 for the new columns I and J
$data = @{
    2  = @(7, 7)
    3  = @(7, 8)
    4  = @(6, 7)
    5  = @(7, 7)
    6  = @(11, 11)
    7  = @(1, 6)
    8  = @(1, 1)
    9  = @(1, 5)
    10 = @(1, 7)
    11 = @(1, 2)
    12 = @(1, 6)
    13 = @(1, 4)
    14 = @(1, 6)
    15 = @(1, 6)
    16 = @(1, 4)
    17 = @(1, 7)
    18 = @(1, 2)
    19 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
